$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 5) of noun-method search results, mirroring the
# existing rows 2-4 in the sheet.
$ws.Range("A5").Value = 42602.583067129628
$ws.Range("B5").Value = "Noun"
$ws.Range("C5").Value = 8928
$ws.Range("D5").Value = 6833
$ws.Range("E5").Value = 1283
$ws.Range("F5").Value = 147
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 66
$ws.Range("M5").Value = 33
